# =========================================================================
# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund holdings detail) positioned right
# before the "总计" (totals) worksheet, and prepends a matching summary row
# to "总计".
# =========================================================================

$wb = $excel.ActiveWorkbook

# Recreate the "总计" sheet from scratch so the newly inserted sheets receive
# sheetId 6 ("2022-Q1") and sheetId 7 ("总计"), matching the target workbook.
$oldTotal = $wb.Worksheets.Item("总计")
[void]$oldTotal.Delete()

$wsQ1 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTotal.Name = "总计"

# -------------------------------------------------------------------------
# Populate "2022-Q1" (fund holdings detail, same shape as other quarters)
# -------------------------------------------------------------------------

# Row 1
$wsQ1.Range("B1").NumberFormat = "@"
$wsQ1.Range("B1").Value = '基金代码'
$wsQ1.Range("B1").ClearFormats()
$wsQ1.Range("C1").NumberFormat = "@"
$wsQ1.Range("C1").Value = '基金名称'
$wsQ1.Range("C1").ClearFormats()
$wsQ1.Range("D1").NumberFormat = "@"
$wsQ1.Range("D1").Value = '基金规模'
$wsQ1.Range("D1").ClearFormats()
$wsQ1.Range("E1").NumberFormat = "@"
$wsQ1.Range("E1").Value = '股票总仓位'
$wsQ1.Range("E1").ClearFormats()
$wsQ1.Range("F1").NumberFormat = "@"
$wsQ1.Range("F1").Value = '仓位占比'
$wsQ1.Range("F1").ClearFormats()
$wsQ1.Range("G1").NumberFormat = "@"
$wsQ1.Range("G1").Value = '持有市值(亿元)'
$wsQ1.Range("G1").ClearFormats()
$wsQ1.Range("H1").NumberFormat = "@"
$wsQ1.Range("H1").Value = '仓位排名'
$wsQ1.Range("H1").ClearFormats()

# Row 2
$wsQ1.Range("A2").Value = 0
$wsQ1.Range("B2").NumberFormat = "@"
$wsQ1.Range("B2").Value = '501054'
$wsQ1.Range("B2").ClearFormats()
$wsQ1.Range("C2").NumberFormat = "@"
$wsQ1.Range("C2").Value = '东方红睿泽三年定期开放灵活配置混合A'
$wsQ1.Range("C2").ClearFormats()
$wsQ1.Range("D2").NumberFormat = "@"
$wsQ1.Range("D2").Value = '109.00'
$wsQ1.Range("D2").ClearFormats()
$wsQ1.Range("E2").NumberFormat = "@"
$wsQ1.Range("E2").Value = '95.90'
$wsQ1.Range("E2").ClearFormats()
$wsQ1.Range("F2").NumberFormat = "@"
$wsQ1.Range("F2").Value = '3.77'
$wsQ1.Range("F2").ClearFormats()
$wsQ1.Range("G2").NumberFormat = "@"
$wsQ1.Range("G2").Value = '4.1093'
$wsQ1.Range("G2").ClearFormats()
$wsQ1.Range("H2").Value = 6

# Row 3
$wsQ1.Range("A3").Value = 1
$wsQ1.Range("B3").NumberFormat = "@"
$wsQ1.Range("B3").Value = '009576'
$wsQ1.Range("B3").ClearFormats()
$wsQ1.Range("C3").NumberFormat = "@"
$wsQ1.Range("C3").Value = '东方红智远三年持有期混合'
$wsQ1.Range("C3").ClearFormats()
$wsQ1.Range("D3").NumberFormat = "@"
$wsQ1.Range("D3").Value = '66.98'
$wsQ1.Range("D3").ClearFormats()
$wsQ1.Range("E3").NumberFormat = "@"
$wsQ1.Range("E3").Value = '92.53'
$wsQ1.Range("E3").ClearFormats()
$wsQ1.Range("F3").NumberFormat = "@"
$wsQ1.Range("F3").Value = '4.76'
$wsQ1.Range("F3").ClearFormats()
$wsQ1.Range("G3").NumberFormat = "@"
$wsQ1.Range("G3").Value = '3.1882'
$wsQ1.Range("G3").ClearFormats()
$wsQ1.Range("H3").Value = 2

# Row 4
$wsQ1.Range("A4").Value = 2
$wsQ1.Range("B4").NumberFormat = "@"
$wsQ1.Range("B4").Value = '005644'
$wsQ1.Range("B4").ClearFormats()
$wsQ1.Range("C4").NumberFormat = "@"
$wsQ1.Range("C4").Value = '广发沪港深行业龙头混合'
$wsQ1.Range("C4").ClearFormats()
$wsQ1.Range("D4").NumberFormat = "@"
$wsQ1.Range("D4").Value = '13.85'
$wsQ1.Range("D4").ClearFormats()
$wsQ1.Range("E4").NumberFormat = "@"
$wsQ1.Range("E4").Value = '88.27'
$wsQ1.Range("E4").ClearFormats()
$wsQ1.Range("F4").NumberFormat = "@"
$wsQ1.Range("F4").Value = '5.38'
$wsQ1.Range("F4").ClearFormats()
$wsQ1.Range("G4").NumberFormat = "@"
$wsQ1.Range("G4").Value = '0.7451'
$wsQ1.Range("G4").ClearFormats()
$wsQ1.Range("H4").Value = 4

# Row 5
$wsQ1.Range("A5").Value = 3
$wsQ1.Range("B5").NumberFormat = "@"
$wsQ1.Range("B5").Value = '008133'
$wsQ1.Range("B5").ClearFormats()
$wsQ1.Range("C5").NumberFormat = "@"
$wsQ1.Range("C5").Value = '华安优质生活混合'
$wsQ1.Range("C5").ClearFormats()
$wsQ1.Range("D5").NumberFormat = "@"
$wsQ1.Range("D5").Value = '10.42'
$wsQ1.Range("D5").ClearFormats()
$wsQ1.Range("E5").NumberFormat = "@"
$wsQ1.Range("E5").Value = '89.58'
$wsQ1.Range("E5").ClearFormats()
$wsQ1.Range("F5").NumberFormat = "@"
$wsQ1.Range("F5").Value = '6.88'
$wsQ1.Range("F5").ClearFormats()
$wsQ1.Range("G5").NumberFormat = "@"
$wsQ1.Range("G5").Value = '0.7169'
$wsQ1.Range("G5").ClearFormats()
$wsQ1.Range("H5").Value = 2

# Row 6
$wsQ1.Range("A6").Value = 4
$wsQ1.Range("B6").NumberFormat = "@"
$wsQ1.Range("B6").Value = '010887'
$wsQ1.Range("B6").ClearFormats()
$wsQ1.Range("C6").NumberFormat = "@"
$wsQ1.Range("C6").Value = '南方消费升级混合A'
$wsQ1.Range("C6").ClearFormats()
$wsQ1.Range("D6").NumberFormat = "@"
$wsQ1.Range("D6").Value = '16.05'
$wsQ1.Range("D6").ClearFormats()
$wsQ1.Range("E6").NumberFormat = "@"
$wsQ1.Range("E6").Value = '78.47'
$wsQ1.Range("E6").ClearFormats()
$wsQ1.Range("F6").NumberFormat = "@"
$wsQ1.Range("F6").Value = '4.34'
$wsQ1.Range("F6").ClearFormats()
$wsQ1.Range("G6").NumberFormat = "@"
$wsQ1.Range("G6").Value = '0.6966'
$wsQ1.Range("G6").ClearFormats()
$wsQ1.Range("H6").Value = 5

# Row 7
$wsQ1.Range("A7").Value = 5
$wsQ1.Range("B7").NumberFormat = "@"
$wsQ1.Range("B7").Value = '006595'
$wsQ1.Range("B7").ClearFormats()
$wsQ1.Range("C7").NumberFormat = "@"
$wsQ1.Range("C7").Value = '广发港股通优质增长混合'
$wsQ1.Range("C7").ClearFormats()
$wsQ1.Range("D7").NumberFormat = "@"
$wsQ1.Range("D7").Value = '8.53'
$wsQ1.Range("D7").ClearFormats()
$wsQ1.Range("E7").NumberFormat = "@"
$wsQ1.Range("E7").Value = '86.63'
$wsQ1.Range("E7").ClearFormats()
$wsQ1.Range("F7").NumberFormat = "@"
$wsQ1.Range("F7").Value = '6.10'
$wsQ1.Range("F7").ClearFormats()
$wsQ1.Range("G7").NumberFormat = "@"
$wsQ1.Range("G7").Value = '0.5203'
$wsQ1.Range("G7").ClearFormats()
$wsQ1.Range("H7").Value = 4

# Row 8
$wsQ1.Range("A8").Value = 6
$wsQ1.Range("B8").NumberFormat = "@"
$wsQ1.Range("B8").Value = '001764'
$wsQ1.Range("B8").ClearFormats()
$wsQ1.Range("C8").NumberFormat = "@"
$wsQ1.Range("C8").Value = '广发沪港深新机遇股票'
$wsQ1.Range("C8").ClearFormats()
$wsQ1.Range("D8").NumberFormat = "@"
$wsQ1.Range("D8").Value = '11.12'
$wsQ1.Range("D8").ClearFormats()
$wsQ1.Range("E8").NumberFormat = "@"
$wsQ1.Range("E8").Value = '92.18'
$wsQ1.Range("E8").ClearFormats()
$wsQ1.Range("F8").NumberFormat = "@"
$wsQ1.Range("F8").Value = '3.93'
$wsQ1.Range("F8").ClearFormats()
$wsQ1.Range("G8").NumberFormat = "@"
$wsQ1.Range("G8").Value = '0.4370'
$wsQ1.Range("G8").ClearFormats()
$wsQ1.Range("H8").Value = 10

# Row 9
$wsQ1.Range("A9").Value = 7
$wsQ1.Range("B9").NumberFormat = "@"
$wsQ1.Range("B9").Value = '910024'
$wsQ1.Range("B9").ClearFormats()
$wsQ1.Range("C9").NumberFormat = "@"
$wsQ1.Range("C9").Value = '东方红启阳三年持有期混合A'
$wsQ1.Range("C9").ClearFormats()
$wsQ1.Range("D9").NumberFormat = "@"
$wsQ1.Range("D9").Value = '6.02'
$wsQ1.Range("D9").ClearFormats()
$wsQ1.Range("E9").NumberFormat = "@"
$wsQ1.Range("E9").Value = '91.72'
$wsQ1.Range("E9").ClearFormats()
$wsQ1.Range("F9").NumberFormat = "@"
$wsQ1.Range("F9").Value = '3.23'
$wsQ1.Range("F9").ClearFormats()
$wsQ1.Range("G9").NumberFormat = "@"
$wsQ1.Range("G9").Value = '0.1944'
$wsQ1.Range("G9").ClearFormats()
$wsQ1.Range("H9").Value = 9

# Row 10
$wsQ1.Range("A10").Value = 8
$wsQ1.Range("B10").NumberFormat = "@"
$wsQ1.Range("B10").Value = '870017'
$wsQ1.Range("B10").ClearFormats()
$wsQ1.Range("C10").NumberFormat = "@"
$wsQ1.Range("C10").Value = '广发资管消费精选灵活配置混合'
$wsQ1.Range("C10").ClearFormats()
$wsQ1.Range("D10").NumberFormat = "@"
$wsQ1.Range("D10").Value = '2.56'
$wsQ1.Range("D10").ClearFormats()
$wsQ1.Range("E10").NumberFormat = "@"
$wsQ1.Range("E10").Value = '93.50'
$wsQ1.Range("E10").ClearFormats()
$wsQ1.Range("F10").NumberFormat = "@"
$wsQ1.Range("F10").Value = '6.42'
$wsQ1.Range("F10").ClearFormats()
$wsQ1.Range("G10").NumberFormat = "@"
$wsQ1.Range("G10").Value = '0.1644'
$wsQ1.Range("G10").ClearFormats()
$wsQ1.Range("H10").Value = 8

# Row 11
$wsQ1.Range("A11").Value = 9
$wsQ1.Range("B11").NumberFormat = "@"
$wsQ1.Range("B11").Value = '010888'
$wsQ1.Range("B11").ClearFormats()
$wsQ1.Range("C11").NumberFormat = "@"
$wsQ1.Range("C11").Value = '南方消费升级混合C'
$wsQ1.Range("C11").ClearFormats()
$wsQ1.Range("D11").NumberFormat = "@"
$wsQ1.Range("D11").Value = '3.39'
$wsQ1.Range("D11").ClearFormats()
$wsQ1.Range("E11").NumberFormat = "@"
$wsQ1.Range("E11").Value = '78.47'
$wsQ1.Range("E11").ClearFormats()
$wsQ1.Range("F11").NumberFormat = "@"
$wsQ1.Range("F11").Value = '4.34'
$wsQ1.Range("F11").ClearFormats()
$wsQ1.Range("G11").NumberFormat = "@"
$wsQ1.Range("G11").Value = '0.1471'
$wsQ1.Range("G11").ClearFormats()
$wsQ1.Range("H11").Value = 5

# Row 12
$wsQ1.Range("A12").Value = 10
$wsQ1.Range("B12").NumberFormat = "@"
$wsQ1.Range("B12").Value = '011032'
$wsQ1.Range("B12").ClearFormats()
$wsQ1.Range("C12").NumberFormat = "@"
$wsQ1.Range("C12").Value = '东方红睿泽三年定期开放灵活配置混合C'
$wsQ1.Range("C12").ClearFormats()
$wsQ1.Range("D12").NumberFormat = "@"
$wsQ1.Range("D12").Value = '0.35'
$wsQ1.Range("D12").ClearFormats()
$wsQ1.Range("E12").NumberFormat = "@"
$wsQ1.Range("E12").Value = '95.90'
$wsQ1.Range("E12").ClearFormats()
$wsQ1.Range("F12").NumberFormat = "@"
$wsQ1.Range("F12").Value = '3.77'
$wsQ1.Range("F12").ClearFormats()
$wsQ1.Range("G12").NumberFormat = "@"
$wsQ1.Range("G12").Value = '0.0132'
$wsQ1.Range("G12").ClearFormats()
$wsQ1.Range("H12").Value = 6

# Row 13
$wsQ1.Range("A13").Value = 11
$wsQ1.Range("B13").NumberFormat = "@"
$wsQ1.Range("B13").Value = '010862'
$wsQ1.Range("B13").ClearFormats()
$wsQ1.Range("C13").NumberFormat = "@"
$wsQ1.Range("C13").Value = '东方红启阳三年持有期混合B'
$wsQ1.Range("C13").ClearFormats()
$wsQ1.Range("E13").NumberFormat = "@"
$wsQ1.Range("E13").Value = '91.72'
$wsQ1.Range("E13").ClearFormats()
$wsQ1.Range("F13").NumberFormat = "@"
$wsQ1.Range("F13").Value = '3.23'
$wsQ1.Range("F13").ClearFormats()
$wsQ1.Range("G13").Value = 0
$wsQ1.Range("H13").Value = 9

# Apply the same cell styling used on the other quarter sheets: bold,
# centered, bordered header row, and bold row-counter column A.
$fmtHeader = $wb.Worksheets.Item("2021-Q4").Range("B1")
$fmtHeader.Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)
$fmtCounter = $wb.Worksheets.Item("2021-Q4").Range("A2")
$fmtCounter.Copy()
$wsQ1.Range("A2:A13").PasteSpecial(-4122)

# -------------------------------------------------------------------------
# Populate "总计" (summary) with the new 2022-Q1 row inserted at the top
# -------------------------------------------------------------------------

# Row 1
$wsTotal.Range("B1").NumberFormat = "@"
$wsTotal.Range("B1").Value = '日期'
$wsTotal.Range("B1").ClearFormats()
$wsTotal.Range("C1").NumberFormat = "@"
$wsTotal.Range("C1").Value = '持有数量(只)'
$wsTotal.Range("C1").ClearFormats()
$wsTotal.Range("D1").NumberFormat = "@"
$wsTotal.Range("D1").Value = '持有市值(亿元)'
$wsTotal.Range("D1").ClearFormats()

# Row 2
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").NumberFormat = "@"
$wsTotal.Range("B2").Value = '2022-Q1'
$wsTotal.Range("B2").ClearFormats()
$wsTotal.Range("C2").Value = 12
$wsTotal.Range("D2").Value = 10.93

# Row 3
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").NumberFormat = "@"
$wsTotal.Range("B3").Value = '2021-Q4'
$wsTotal.Range("B3").ClearFormats()
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 3.37

# Row 4
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").NumberFormat = "@"
$wsTotal.Range("B4").Value = '2021-Q3'
$wsTotal.Range("B4").ClearFormats()
$wsTotal.Range("C4").Value = 14
$wsTotal.Range("D4").Value = 14.16

# Row 5
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").NumberFormat = "@"
$wsTotal.Range("B5").Value = '2021-Q2'
$wsTotal.Range("B5").ClearFormats()
$wsTotal.Range("C5").Value = 17
$wsTotal.Range("D5").Value = 14.1

# Row 6
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").NumberFormat = "@"
$wsTotal.Range("B6").Value = '2021-Q1'
$wsTotal.Range("B6").ClearFormats()
$wsTotal.Range("C6").Value = 14
$wsTotal.Range("D6").Value = 11.71

# Row 7
$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").NumberFormat = "@"
$wsTotal.Range("B7").Value = '2020-Q4'
$wsTotal.Range("B7").ClearFormats()
$wsTotal.Range("C7").Value = 8
$wsTotal.Range("D7").Value = 3.44

$fmtHeaderTotal = $wb.Worksheets.Item("2021-Q4").Range("B1")
$fmtHeaderTotal.Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$fmtCounterTotal = $wb.Worksheets.Item("2021-Q4").Range("A2")
$fmtCounterTotal.Copy()
$wsTotal.Range("A2:A7").PasteSpecial(-4122)

# Select A1 on the new detail sheet to mirror the other sheets' default view
[void]$wsQ1.Range("A1").Select()

